$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Activate()

$summary.Range("A3").Value = 301.72000000000003
$summary.Range("E3").Value = 301.72000000000003

# Update the active selection on the Summary sheet
$summary.Range("C10").Select()

# ---- Repayment schedule sheet ----
$sched = $wb.Worksheets.Item("Repayment schedule")
$sched.Activate()

# Row 11
$sched.Range("B11").Value = 14
$sched.Range("C11").Value = 42157
$sched.Range("F11").Value = 842.84
$sched.Range("G11").Value = 2551.96
$sched.Range("H11").Value = 15.63

# Row 12
$sched.Range("B12").Value = 14
$sched.Range("F12").Value = 846.72
$sched.Range("G12").Value = 1705.24
$sched.Range("H12").Value = 11.75

# Row 13
$sched.Range("F13").Value = 850.62
$sched.Range("G13").Value = 854.62
$sched.Range("H13").Value = 7.85

# Row 14
$sched.Range("F14").Value = 854.62
$sched.Range("H14").Value = 3.93
$sched.Range("K14").Value = 858.55
$sched.Range("Q14").Value = 858.55

# Update the active selection on the Repayment schedule sheet
$sched.Range("J15").Select()
